$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 2869.9644
$ws.Range("I45").Value = 2921.913
$ws.Range("J45").Value = 2631
$ws.Range("K45").Value = 2921.913
$ws.Range("L45").Value = 2631
$ws.Range("M45").Value = -2544.913
$ws.Range("N45").Value = -3385
$ws.Range("H122").Value = 5052594
$ws.Range("I122").Value = 2236.1667
$ws.Range("J122").Value = 27779204
$ws.Range("K122").Value = 6708.500100000001
$ws.Range("L122").Value = 83337612
$ws.Range("M122").Value = -4258.500100000001
$ws.Range("N122").Value = -83342512
$ws.Range("H139").Value = 65056
$ws.Range("J139").Value = 65056
$ws.Range("L139").Value = 65056
$ws.Range("N139").Value = -75336

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 1114.5834
$ws.Range("I20").Value = 1055.7
$ws.Range("J20").Value = 1409
$ws.Range("K20").Value = 1055.7
$ws.Range("L20").Value = 1409
$ws.Range("M20").Value = -808.7
$ws.Range("N20").Value = -1903
$ws.Range("H132").Value = 44026.668
$ws.Range("J132").Value = 44026.668
$ws.Range("L132").Value = 44026.668
$ws.Range("N132").Value = -54146.668

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H20").Value = 51000
$ws.Range("J20").Value = 51000
$ws.Range("L20").Value = 51000
$ws.Range("N20").Value = -51472
$ws.Range("H30").Value = 51000
$ws.Range("J30").Value = 51000
$ws.Range("L30").Value = 51000
$ws.Range("N30").Value = -51182
$ws.Range("H122").Value = 1264.0938
$ws.Range("I122").Value = 1175.8387
$ws.Range("K122").Value = 3527.5161
$ws.Range("M122").Value = -1077.5161
$ws.Range("H127").Value = 33181.816
$ws.Range("J127").Value = 33181.816
$ws.Range("L127").Value = 33181.816
$ws.Range("N127").Value = -43101.816
$ws.Range("H128").Value = 51000
$ws.Range("J128").Value = 51000
$ws.Range("L128").Value = 51000
$ws.Range("N128").Value = -60960
$ws.Range("H135").Value = 59775
$ws.Range("J135").Value = 59775
$ws.Range("L135").Value = 59775
$ws.Range("N135").Value = -69915
$ws.Range("H138").Value = 49773.332
$ws.Range("J138").Value = 49773.332
$ws.Range("L138").Value = 49773.332
$ws.Range("N138").Value = -60053.332
$ws.Range("H140").Value = 62778.57
$ws.Range("J140").Value = 62778.57
$ws.Range("L140").Value = 62778.57
$ws.Range("N140").Value = -73138.57000000001

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 19023.926
$ws.Range("I5").Value = 29769.47
$ws.Range("J5").Value = 756.5
$ws.Range("K5").Value = 89308.41
$ws.Range("L5").Value = 2269.5
$ws.Range("M5").Value = -89196.41
$ws.Range("N5").Value = -2493.5
$ws.Range("H96").Value = 4250
$ws.Range("J96").Value = 4250
$ws.Range("L96").Value = 12750
$ws.Range("N96").Value = -16868
$ws.Range("H105").Value = 12000
$ws.Range("I105").Value = 0
$ws.Range("J105").Value = 12000
$ws.Range("K105").Value = 0
$ws.Range("L105").Value = 36000
$ws.Range("M105").Value = ""
$ws.Range("N105").Value = -41242
$ws.Range("H113").Value = 639.45
$ws.Range("I113").Value = 600
$ws.Range("J113").Value = 641.5263
$ws.Range("K113").Value = 1800
$ws.Range("L113").Value = 1924.5789
$ws.Range("M113").Value = 370
$ws.Range("N113").Value = -6264.5789
$ws.Range("H131").Value = 683.7222
$ws.Range("I131").Value = 408.8889
$ws.Range("J131").Value = 958.55554
$ws.Range("K131").Value = 1226.6667
$ws.Range("L131").Value = 2875.66662
$ws.Range("M131").Value = 3813.3333
$ws.Range("N131").Value = -12955.66662
$ws.Range("H134").Value = 1852.4375
$ws.Range("I134").Value = 1713.7333
$ws.Range("J134").Value = 3933
$ws.Range("K134").Value = 5141.199900000001
$ws.Range("L134").Value = 11799
$ws.Range("M134").Value = -71.19990000000053
$ws.Range("N134").Value = -21939
$ws.Range("H135").Value = 19023.926
$ws.Range("I135").Value = 29769.47
$ws.Range("J135").Value = 756.5
$ws.Range("K135").Value = 267925.23
$ws.Range("L135").Value = 6808.5
$ws.Range("M135").Value = -265390.23
$ws.Range("N135").Value = -11878.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 2398.8965
$ws.Range("I122").Value = 2029.579
$ws.Range("K122").Value = 6088.737
$ws.Range("M122").Value = -3638.737
$ws.Range("H133").Value = 63181.91
$ws.Range("J133").Value = 63181.91
$ws.Range("L133").Value = 63181.91
$ws.Range("N133").Value = -73301.91

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 474.5625
$ws.Range("I22").Value = 244.88889
$ws.Range("J22").Value = 769.8570999999999
$ws.Range("K22").Value = 244.88889
$ws.Range("L22").Value = 769.8570999999999
$ws.Range("M22").Value = 50.11111
$ws.Range("N22").Value = -1359.8571
$ws.Range("H27").Value = 474.5625
$ws.Range("I27").Value = 244.88889
$ws.Range("J27").Value = 769.8570999999999
$ws.Range("K27").Value = 244.88889
$ws.Range("L27").Value = 769.8570999999999
$ws.Range("M27").Value = -137.88889
$ws.Range("N27").Value = -983.8570999999999
$ws.Range("H40").Value = 5666.3335
$ws.Range("I40").Value = 6166.5
$ws.Range("K40").Value = 6166.5
$ws.Range("M40").Value = -6030.5
$ws.Range("H132").Value = 46474.434
$ws.Range("I132").Value = 1956.7778
$ws.Range("K132").Value = 5870.3334
$ws.Range("M132").Value = -3340.3334

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 414.81818
$ws.Range("I107").Value = 382.875
$ws.Range("J107").Value = 500
$ws.Range("K107").Value = 1148.625
$ws.Range("L107").Value = 1500
$ws.Range("M107").Value = 771.375
$ws.Range("N107").Value = -5340
$ws.Range("H122").Value = 3664.6843
$ws.Range("I122").Value = 3425
$ws.Range("J122").Value = 3880.4
$ws.Range("K122").Value = 10275
$ws.Range("L122").Value = 11641.2
$ws.Range("M122").Value = -7825
$ws.Range("N122").Value = -16541.2
$ws.Range("H126").Value = 2099.75
$ws.Range("I126").Value = 1242.4286
$ws.Range("J126").Value = 3300
$ws.Range("K126").Value = 3727.2858
$ws.Range("L126").Value = 9900
$ws.Range("M126").Value = -1257.2858
$ws.Range("N126").Value = -14840
$ws.Range("H132").Value = 66550.64999999999
$ws.Range("I132").Value = 38950.63
$ws.Range("K132").Value = 116851.89
$ws.Range("M132").Value = -114321.89
